$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Change 1: "Definisci un set minimo di dati anagrafici " paragraph.
# The run " un set " is split into " " + "un set" + " " (the middle
# piece gets wrapped by the grammar-checker's gramStart/gramEnd marks
# in real Word; here we reproduce the run split itself, which is the
# structural, content-visible part of the edit).
# --------------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute(" un set ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $outerStart = $find1.Start
    $outerEnd = $find1.End
    # Sub-range covering just "un set" (without the surrounding spaces)
    $mid = $d.Range($outerStart + 1, $outerEnd - 1)
    # Toggling Bold off then back on forces Word to split the run at
    # these boundaries while leaving the resulting formatting exactly
    # as it was (Bold stays on, nothing else changes).
    $mid.Bold = 0
    $mid.Bold = 1
}

# --------------------------------------------------------------------
# Change 2: "Questionario minimo per la raccolta dei dati anagrafici
# dei partecipanti al test:" paragraph drops from 14pt (sz 28) to
# 13pt (sz 26) and gains an explicit complex-script size (szCs 26),
# including on the paragraph mark itself. Locate the paragraph by
# scanning (rather than Find's own .Paragraphs collection) so that the
# resulting Range includes the trailing paragraph mark and the size
# change lands on the paragraph mark's run properties too.
# --------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Questionario minimo per la raccolta dei dati anagrafici dei partecipanti al test:*") {
        $pRng = $p.Range
        $pRng.Font.Size = 13
        $pRng.Font.SizeBi = 13
        break
    }
}
